$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201, shifting existing rows 201-272 down to 202-273.
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new record.
$ws.Range("A201").Value = 9
$ws.Range("B201").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C201").Value = "Metropolitana"
$ws.Range("D201").Value = 44559
$ws.Range("E201").Value = 13
$ws.Range("F201").Value = 100112052
$ws.Range("G201").Value = "Albahaca"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 160
$ws.Range("K201").Value = 4000
$ws.Range("L201").Value = 4500
$ws.Range("M201").Value = 4250
$ws.Range("N201").Value = "`$/docena de matas"
$ws.Range("O201").Value = "Región Metropolitana"
$ws.Range("P201").Value = 708
$ws.Range("Q201").Value = 6
$ws.Range("R201").Value = "Hortaliza"
